$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Terry Rozier", "PG", "Miami Heat"),
    @("Donte DiVincenzo", "SG,SF", "Minnesota Timberwolves"),
    @("Dalton Knecht", "SG,SF", "Los Angeles Lakers"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Norman Powell", "SG,SF", "LA Clippers")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
